$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 (Ватаманюк Анна): mark C9:G9 as "ок" (H9 already was "ок")
$ws.Range("C9:G9").Value = "ок"

# Row 16 (Конова Елизавета): mark E16:G16 as "ок" (H16 already was "ок")
$ws.Range("E16:G16").Value = "ок"

# Update the active selection to H9 as recorded in the workbook view
$ws.Range("H9").Select()
